$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 (row label "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 243
$wsOff.Range("C3").Value = 181
$wsOff.Range("D3").Value = 59
$wsOff.Range("E3").Value = 32
$wsOff.Range("F3").Value = 5

# DEF sheet - row 3 (row label "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 254
$wsDef.Range("C3").Value = 183
$wsDef.Range("D3").Value = 47
$wsDef.Range("E3").Value = 26
